# "fix extra time sum" - the totals row ("مجموع" / sum row) was sitting two
# rows below the header, separated by blank rows 2 and 3. Remove one of the
# blank rows above it so the totals row shifts up from row 4 to row 3,
# directly under the header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Delete()

